# Update the Jogos do Dia Betfair Back/Lay workbook:
# - Row 2 now represents what used to be the second match (Fortaleza FC x
#   Alianza FC Valledupar) with freshly scraped odds.
# - The old row 3 (which held that match previously) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 with the new match info / odds ---
$ws.Range("C2").Value = "22:30:00"
$ws.Range("D2").Value = "Fortaleza FC"
$ws.Range("E2").Value = "Alianza FC Valledupar"

$ws.Range("F2").Value = 2.46
$ws.Range("G2").Value = 2.82
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 6.8
$ws.Range("J2").Value = 2.14
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 1.43
$ws.Range("N2").Value = 1.5
$ws.Range("O2").Value = 2.94
$ws.Range("P2").Value = 1.12
$ws.Range("Q2").Value = 8.6
$ws.Range("R2").Value = 1.02
$ws.Range("S2").Value = 32
$ws.Range("T2").Value = 4.8
$ws.Range("U2").Value = 1.2
$ws.Range("V2").Value = 1.2
$ws.Range("W2").Value = 1.64
$ws.Range("X2").Value = 3.35
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 4.8
$ws.Range("AC2").Value = 9.4
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 12.5
$ws.Range("AG2").Value = 950
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# --- Remove the old row 3 (data now consolidated into row 2) ---
$ws.Rows("3:3").Delete()
